$d = $word.ActiveDocument

function Set-ParagraphXml {
    # Rebuilds a paragraph's content in-place from an explicit OOXML <w:p> ...
    # </w:p> fragment (preserving <w:pPr> / run formatting exactly as given),
    # so that zero-length empty runs (<w:r/>) that a plain Find/Replace would
    # otherwise silently merge away stay intact.
    param($doc, $paraIndex, $innerP)

    $p = $doc.Paragraphs.Item($paraIndex)
    $rOrig = $p.Range
    $isLast = ($paraIndex -eq $doc.Paragraphs.Count)

    if ($isLast) {
        # The very last paragraph's Range also swallows the final paragraph
        # mark; excluding it keeps the paragraph count/structure unchanged
        # (otherwise InsertXML would leave a stray empty paragraph behind).
        $r = $doc.Range($rOrig.Start, $rOrig.End - 1)
    } else {
        $r = $doc.Range($rOrig.Start, $rOrig.End)
    }

    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerP + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# 1. Heading1 title (no surrounding empty runs, so a plain Find/Replace keeps
#    the structure intact)
$d.Content.Find.Execute(
    "Play Cash Compass for Free - Exciting Pirate-Themed Slot Game",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Cash Compass for Free", 1)

# 2. "What we like" bullets (paragraphs contain a leading empty <w:r/> run
#    that must be preserved, so rebuild each paragraph verbatim)
Set-ParagraphXml $d 42 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Inspired by pirates and deserted islands of the Caribbean</w:t></w:r></w:p>'

Set-ParagraphXml $d 44 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Various bonus games available</w:t></w:r></w:p>'

Set-ParagraphXml $d 45 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Maximum win of over 7,400 times the bet</w:t></w:r></w:p>'

# 3. "What we don't like" bullets
Set-ParagraphXml $d 47 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Limited betting options</w:t></w:r></w:p>'

Set-ParagraphXml $d 48 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Limited autoplay settings</w:t></w:r></w:p>'

# 4. Bold title line near the end (leading empty <w:r/> run + bold run)
Set-ParagraphXml $d 49 '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cash Compass for Free</w:t></w:r></w:p>'

# 5. Italic meta-description line (last paragraph in the body - leading empty
#    <w:r/> run is reconstructed automatically, see Set-ParagraphXml above)
Set-ParagraphXml $d 50 '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Cash Compass, a slot game inspired by pirates and deserted islands. Play for free now!</w:t></w:r></w:p>'
